# Latest update 19-04-2023 eve
# Adds a new "Add_Products" worksheet and appends new user/partner/product records
# to the existing Add_User and Login worksheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Add_User
$ws2 = $wb.Worksheets.Item(2)   # Login

# Create the new "Add_Products" worksheet after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "Add_Products"

$ws3.Range("A1").Value = 'Product_Identifier_Value'
$ws3.Range("A2").Value = 'ProdIdoSu1'
$ws3.Range("A3").Value = 'ProdIdhpsv'
$ws3.Range("A4").Value = 'ProdIdglSB'
$ws3.Range("A5").Value = 'ProdIdkU8U'
$ws3.Range("A6").Value = 'ProdIdJnJC'
$ws3.Range("A7").Value = 'ProdIdJzWT'
$ws3.Range("A8").Value = 'ProdIdonIg'
$ws3.Range("A9").Value = 'ProdIdvcWQ'
$ws3.Range("A10").Value = 'ProdIdARpH'
$ws3.Range("A11").Value = 'ProdIdsoSn'
$ws3.Range("A12").Value = 'ProdIdPjsl'
$ws1.Range("C12").Value = 'digitalmeshRtUaWD@maildrop.cc'
$ws2.Range("D58").Value = 'digitalmeshfdykxm@maildrop.cc'
$ws3.Range("A13").Value = 'ProdIdlN8Y'
$ws2.Range("F18").NumberFormat = "@"
$ws2.Range("F18").Value = '1196203914'
$ws2.Range("F18").ClearFormats()
$ws2.Range("C79").Value = 'jAuJT'
$ws2.Range("D59").Value = 'digitalmesh7tvwde@maildrop.cc'
$ws2.Range("F19").NumberFormat = "@"
$ws2.Range("F19").Value = '8243655906'
$ws2.Range("F19").ClearFormats()
$ws2.Range("C80").Value = 'ebGjf'
$ws2.Range("D60").Value = 'digitalmeshezwzur@maildrop.cc'
$ws2.Range("F20").NumberFormat = "@"
$ws2.Range("F20").Value = '9628248577'
$ws2.Range("F20").ClearFormats()
$ws2.Range("C81").Value = 'vCtJo'
$ws2.Range("D61").Value = 'digitalmeshctmjob@maildrop.cc'
$ws2.Range("F21").NumberFormat = "@"
$ws2.Range("F21").Value = '5536733521'
$ws2.Range("F21").ClearFormats()
$ws2.Range("C82").Value = 'DKaVL'
$ws2.Range("D62").Value = 'digitalmeshangeb7@maildrop.cc'
$ws2.Range("F22").NumberFormat = "@"
$ws2.Range("F22").Value = '2622922229'
$ws2.Range("F22").ClearFormats()
$ws2.Range("C83").Value = 'UGCVF'
$ws2.Range("E17").Value = 'DmBzMHPs!2'
$ws3.Range("A14").Value = 'ProdIdUuYv'
$ws3.Range("A15").Value = 'ProdIdBxgT'
$ws3.Range("A16").Value = 'ProdIdBkhn'
$ws3.Range("A17").Value = 'ProdIdzzq3'
$ws3.Range("A18").Value = 'ProdIdJdAy'
$ws3.Range("A19").Value = 'ProdIdbvkb'
$ws3.Range("A20").Value = 'ProdIdQWs0'
$ws3.Range("A21").Value = 'ProdIdt01S'
$ws3.Range("A22").Value = 'ProdIdfCFB'
$ws3.Range("A23").Value = 'ProdId7yCE'
$ws3.Range("A24").Value = 'ProdIdGpTv'
$ws3.Range("A25").Value = 'ProdIdbNMQ'
$ws3.Range("A26").Value = 'ProdIdAbiM'
$ws3.Range("A27").Value = 'ProdIdQ06U'
$ws3.Range("A28").Value = 'ProdIdcjtX'
$ws3.Range("A29").Value = 'ProdIdBI5b'
$ws3.Range("A30").Value = 'ProdIdGcPx'
# Approximate the original author's column auto-fit sizing for the new sheet.
$ws3.Columns.Item(1).EntireColumn.AutoFit() | Out-Null
